$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = 251752
$ws.Range("B4").Value = "BIMEC 5"
$ws.Range("C4").Value = 17
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = "2025-06-04 12:00:00"
$ws.Range("F4").Value = "2025-06-04 12:17:00"
$ws.Range("G4").Value = "2025-06-04 12:17:00"
$ws.Range("H4").Value = "2025-06-04 12:17:00"
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = "bobina"
$ws.Range("K4").Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;R12 ;R6 ;R9"
$ws.Range("L4").Value = 3
$ws.Range("M4").Value = 76
$ws.Range("N4").Value = 39846
$ws.Range("O4").Value = "X"
$ws.Range("P4").Value = 39846
$ws.Range("Q4").Value = "2025-05-20 00:00:00"
$ws.Range("R4").Value = -0.5118055555555555
$ws.Range("S4").Value = 1

# Row 5
$ws.Range("A5").Value = 251218
$ws.Range("B5").Value = "BIMEC 5"
$ws.Range("C5").Value = 21
$ws.Range("D5").Value = 96.90140845070422
$ws.Range("E5").Value = "2025-06-04 12:17:00"
$ws.Range("F5").Value = "2025-06-04 12:38:00"
$ws.Range("G5").Value = "2025-06-04 12:38:00"
$ws.Range("H5").Value = "2025-06-04 14:14:54"
$ws.Range("I5").Value = 6880
$ws.Range("J5").Value = "bobina"
$ws.Range("K5").Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;R12 ;R9"
$ws.Range("L5").Value = 6
$ws.Range("M5").Value = 76
$ws.Range("N5").Value = 39885
$ws.Range("O5").Value = "X"
$ws.Range("P5").Value = 39885
$ws.Range("Q5").Value = "2025-05-09 00:00:00"
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 1

# Row 6
$ws.Range("A6").Value = 251895
$ws.Range("B6").Value = "BIMEC 5"
$ws.Range("C6").Value = 38
$ws.Range("D6").Value = 249.2112676056338
$ws.Range("E6").Value = "2025-06-04 14:14:54"
$ws.Range("F6").Value = "2025-06-04 14:52:54"
$ws.Range("G6").Value = "2025-06-04 14:52:54"
$ws.Range("H6").Value = "2025-06-05 11:02:06"
$ws.Range("I6").Value = 17694
$ws.Range("J6").Value = "bobina"
$ws.Range("K6").Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R6 ;R9"
$ws.Range("L6").Value = 10
$ws.Range("M6").Value = 76
$ws.Range("N6").Value = "39891 (esterno)"
$ws.Range("O6").Value = "X"
$ws.Range("P6").Value = 39891
$ws.Range("Q6").Value = "2025-05-26 00:00:00"
$ws.Range("R6").Value = -10.45980046947917
$ws.Range("S6").Value = 4

# Row 7
$ws.Range("A7").Value = 251565
$ws.Range("B7").Value = "BIMEC 5"
$ws.Range("C7").Value = 31
$ws.Range("D7").Value = 176.7464788732394
$ws.Range("E7").Value = "2025-06-05 11:02:06"
$ws.Range("F7").Value = "2025-06-05 11:33:06"
$ws.Range("G7").Value = "2025-06-05 11:33:06"
$ws.Range("H7").Value = "2025-06-05 14:29:51"
$ws.Range("I7").Value = 12549
$ws.Range("J7").Value = "bobina"
$ws.Range("K7").Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$ws.Range("L7").Value = 2
$ws.Range("M7").Value = 70
$ws.Range("N7").Value = 39885
$ws.Range("O7").Value = "X"
$ws.Range("P7").Value = 39885
$ws.Range("Q7").Value = "2025-06-10 00:00:00"
$ws.Range("R7").Value = -0.6040688575925925
$ws.Range("S7").Value = 1

# Row 8
$ws.Range("A8").Value = 251362
$ws.Range("B8").Value = "BIMEC 5"
$ws.Range("C8").Value = 32
$ws.Range("D8").Value = 35.28169014084507
$ws.Range("E8").Value = "2025-06-05 14:29:51"
$ws.Range("F8").Value = "2025-06-06 07:01:51"
$ws.Range("G8").Value = "2025-06-06 07:01:51"
$ws.Range("H8").Value = "2025-06-06 07:37:08"
$ws.Range("I8").Value = 2505
$ws.Range("J8").Value = "bobina"
$ws.Range("K8").Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;R12 ;R9"
$ws.Range("L8").Value = 3
$ws.Range("M8").Value = 76
$ws.Range("N8").Value = 39874
$ws.Range("O8").Value = "X"
$ws.Range("P8").Value = 39874
$ws.Range("Q8").Value = "2025-04-24 00:00:00"
$ws.Range("R8").Value = -1.317458920185185
$ws.Range("S8").Value = 7

# Row 9
$ws.Range("A9").Value = 251780
$ws.Range("B9").Value = "BIMEC 5"
$ws.Range("C9").Value = 30
$ws.Range("D9").Value = 342.2394366197183
$ws.Range("E9").Value = "2025-06-06 07:37:08"
$ws.Range("F9").Value = "2025-06-06 08:07:08"
$ws.Range("G9").Value = "2025-06-06 08:07:08"
$ws.Range("H9").Value = "2025-06-06 13:49:22"
$ws.Range("I9").Value = 24299
$ws.Range("J9").Value = "bobina"
$ws.Range("K9").Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$ws.Range("L9").Value = 3
$ws.Range("M9").Value = 70
$ws.Range("N9").Value = "39887 (esterno)"
$ws.Range("O9").Value = "X"
$ws.Range("P9").Value = 39887
$ws.Range("Q9").Value = "2025-06-18 00:00:00"
$ws.Range("R9").Value = 0
$ws.Range("S9").Value = 7

# Row 10
$ws.Range("A10").Value = 251070
$ws.Range("B10").Value = "CASON"
$ws.Range("C10").Value = 36.5
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = "2025-06-05 07:00:00"
$ws.Range("F10").Value = "2025-06-05 07:36:30"
$ws.Range("G10").Value = "2025-06-05 07:36:30"
$ws.Range("H10").Value = "2025-06-05 07:36:30"
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = "bobina"
$ws.Range("K10").Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R6 ;R9"
$ws.Range("L10").Value = 6
$ws.Range("M10").Value = 70
$ws.Range("N10").Value = 39885
$ws.Range("O10").Value = "X"
$ws.Range("P10").Value = 39885
$ws.Range("Q10").Value = "2025-03-28 00:00:00"
$ws.Range("R10").Value = -0.3170138888888889
$ws.Range("S10").Value = 2

# Row 11
$ws.Range("A11").Value = 251773
$ws.Range("B11").Value = "CASON"
$ws.Range("C11").Value = 32.5
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = "2025-06-05 07:36:30"
$ws.Range("F11").Value = "2025-06-05 08:09:00"
$ws.Range("G11").Value = "2025-06-05 08:09:00"
$ws.Range("H11").Value = "2025-06-05 08:09:00"
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = "bobina"
$ws.Range("K11").Value = "CASON ;R6"
$ws.Range("L11").Value = 7
$ws.Range("M11").Value = 70
$ws.Range("N11").Value = 39874
$ws.Range("O11").Value = "X"
$ws.Range("P11").Value = 39874
$ws.Range("Q11").Value = "2025-05-25 00:00:00"
$ws.Range("R11").Value = -0.3395833333333333
$ws.Range("S11").Value = 1

# Row 12
$ws.Range("A12").Value = 251500
$ws.Range("B12").Value = "CASON"
$ws.Range("C12").Value = 36.5
$ws.Range("D12").Value = 179.9272727272727
$ws.Range("E12").Value = "2025-06-05 08:09:00"
$ws.Range("F12").Value = "2025-06-05 08:45:30"
$ws.Range("G12").Value = "2025-06-05 08:45:30"
$ws.Range("H12").Value = "2025-06-05 11:45:25"
$ws.Range("I12").Value = 9896
$ws.Range("J12").Value = "bobina"
$ws.Range("K12").Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R6 ;R9"
$ws.Range("L12").Value = 4
$ws.Range("M12").Value = 70
$ws.Range("N12").Value = 39885
$ws.Range("O12").Value = "X"
$ws.Range("P12").Value = 39885
$ws.Range("Q12").Value = "2025-05-26 00:00:00"
$ws.Range("R12").Value = -0.4898800505092593
$ws.Range("S12").Value = 2

# Row 13
$ws.Range("A13").Value = 251180
$ws.Range("B13").Value = "R10"
$ws.Range("C13").Value = 30
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = "2025-06-04 07:00:00"
$ws.Range("F13").Value = "2025-06-04 07:30:00"
$ws.Range("G13").Value = "2025-06-04 07:30:00"
$ws.Range("H13").Value = "2025-06-04 07:30:00"
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = "bobina"
$ws.Range("K13").Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$ws.Range("L13").Value = 4
$ws.Range("M13").Value = 70
$ws.Range("N13").Value = "39887 (esterno)"
$ws.Range("O13").Value = "X"
$ws.Range("P13").Value = 39887
$ws.Range("Q13").Value = "2025-05-20 00:00:00"
$ws.Range("R13").Value = -15.3125
$ws.Range("S13").Value = 7

# Row 14
$ws.Range("A14").Value = 252282
$ws.Range("B14").Value = "R10"
$ws.Range("C14").Value = 25
$ws.Range("D14").Value = 44.88524590163934
$ws.Range("E14").Value = "2025-06-04 07:30:00"
$ws.Range("F14").Value = "2025-06-04 07:55:00"
$ws.Range("G14").Value = "2025-06-04 07:55:00"
$ws.Range("H14").Value = "2025-06-04 08:39:53"
$ws.Range("I14").Value = 2738
$ws.Range("J14").Value = "bobina"
$ws.Range("K14").Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$ws.Range("L14").Value = 5
$ws.Range("M14").Value = 70
$ws.Range("N14").Value = 39885
$ws.Range("O14").Value = "X"
$ws.Range("P14").Value = 39885
$ws.Range("Q14").Value = "2025-06-09 00:00:00"
$ws.Range("R14").Value = 0
$ws.Range("S14").Value = 1

# Row 15
$ws.Range("A15").Value = 251984
$ws.Range("B15").Value = "R10"
$ws.Range("C15").Value = 30
$ws.Range("D15").Value = 338.327868852459
$ws.Range("E15").Value = "2025-06-04 08:39:53"
$ws.Range("F15").Value = "2025-06-04 09:09:53"
$ws.Range("G15").Value = "2025-06-04 09:09:53"
$ws.Range("H15").Value = "2025-06-04 14:48:12"
$ws.Range("I15").Value = 20638
$ws.Range("J15").Value = "bobina"
$ws.Range("K15").Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$ws.Range("L15").Value = 3
$ws.Range("M15").Value = 70
$ws.Range("N15").Value = 39874
$ws.Range("O15").Value = "X"
$ws.Range("P15").Value = 39874
$ws.Range("Q15").Value = "2025-06-10 00:00:00"
$ws.Range("R15").Value = 0
$ws.Range("S15").Value = 1

# Row 16
$ws.Range("A16").Value = 252084
$ws.Range("B16").Value = "R10"
$ws.Range("C16").Value = 25
$ws.Range("D16").Value = 641
$ws.Range("E16").Value = "2025-06-04 14:48:12"
$ws.Range("F16").Value = "2025-06-05 07:13:12"
$ws.Range("G16").Value = "2025-06-05 07:13:12"
$ws.Range("H16").Value = "2025-06-06 09:54:12"
$ws.Range("I16").Value = 39101
$ws.Range("J16").Value = "bobina"
$ws.Range("K16").Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$ws.Range("L16").Value = 2
$ws.Range("M16").Value = 70
$ws.Range("N16").Value = 39885
$ws.Range("O16").Value = "X"
$ws.Range("P16").Value = 39885
$ws.Range("Q16").Value = "2025-06-30 00:00:00"
$ws.Range("R16").Value = -1.412647996354166
$ws.Range("S16").Value = 7
